$wb = $excel.ActiveWorkbook

# --- Sheet "data": update Профвзносы (B), Месяц (C), Сумма (D) ---
$wsData = $wb.Worksheets("data")

$wsData.Range("B2").Value = "+"
$wsData.Range("C2").Value = "Все"
$wsData.Range("D2").Value = 430
$wsData.Range("C3").Value = "Апрель"
$wsData.Range("D3").Value = 278
$wsData.Range("B4").Value = "-"
$wsData.Range("C4").Value = "Декабрь"
$wsData.Range("D4").Value = 223
$wsData.Range("D5").Value = 385
$wsData.Range("B6").Value = "-"
$wsData.Range("C6").Value = "Декабрь"
$wsData.Range("D6").Value = 229
$wsData.Range("C7").Value = "Апрель"
$wsData.Range("D7").Value = 292
$wsData.Range("B8").Value = "-"
$wsData.Range("C8").Value = "Декабрь"
$wsData.Range("D8").Value = 295
$wsData.Range("C9").Value = "Все"
$wsData.Range("D9").Value = 82
$wsData.Range("C10").Value = "Все"
$wsData.Range("D10").Value = 100
$wsData.Range("C11").Value = "Февраль"
$wsData.Range("D11").Value = 357
$wsData.Range("B12").Value = "-"
$wsData.Range("C12").Value = "Сентябрь"
$wsData.Range("D12").Value = 313
$wsData.Range("C13").Value = "Все"
$wsData.Range("D13").Value = 119
$wsData.Range("D14").Value = 264
$wsData.Range("C15").Value = "Декабрь"
$wsData.Range("D15").Value = 225
$wsData.Range("B16").Value = "+"
$wsData.Range("D16").Value = 250
$wsData.Range("B17").Value = "-"
$wsData.Range("C17").Value = "Декабрь"
$wsData.Range("D17").Value = 84
$wsData.Range("B18").Value = "-"
$wsData.Range("C18").Value = "Сентябрь"
$wsData.Range("D18").Value = 301
$wsData.Range("C19").Value = "Апрель"
$wsData.Range("D19").Value = 442
$wsData.Range("C20").Value = "Апрель"
$wsData.Range("D20").Value = 309
$wsData.Range("C21").Value = "Апрель"
$wsData.Range("D21").Value = 316
$wsData.Range("C22").Value = "Апрель"
$wsData.Range("D22").Value = 286
$wsData.Range("B23").Value = "-"
$wsData.Range("C23").Value = "Сентябрь"
$wsData.Range("D23").Value = 295
$wsData.Range("C24").Value = "Февраль"
$wsData.Range("D24").Value = 292
$wsData.Range("D25").Value = 253
$wsData.Range("C26").Value = "Февраль"
$wsData.Range("D26").Value = 339
$wsData.Range("B27").Value = "-"
$wsData.Range("C27").Value = "Апрель"
$wsData.Range("D27").Value = 412
$wsData.Range("C28").Value = "Февраль"
$wsData.Range("D28").Value = 301
$wsData.Range("C29").Value = "Сентябрь"
$wsData.Range("D29").Value = 298
$wsData.Range("C30").Value = "Декабрь"
$wsData.Range("D30").Value = 247

# --- Sheet "results": update Санаторий (B), Статус (C) ---
$wsResults = $wb.Worksheets("results")

$wsResults.Range("B2").Value = "Лесной"
$wsResults.Range("C2").Value = "едет"
$wsResults.Range("B4").Value = "Искра"
$wsResults.Range("C4").Value = "не едет"
$wsResults.Range("B6").Value = "Янтарь"
$wsResults.Range("C6").Value = "не едет"
$wsResults.Range("B7").Value = "Янтарь"
$wsResults.Range("B8").Value = "Искра"
$wsResults.Range("C8").Value = "не едет"
$wsResults.Range("B9").Value = "Искра"
$wsResults.Range("B10").Value = "Прибой"
$wsResults.Range("B11").Value = "Искра"
$wsResults.Range("B12").Value = "Янтарь"
$wsResults.Range("C12").Value = "не едет"
$wsResults.Range("B13").Value = "Лесной"
$wsResults.Range("B15").Value = "Волна"
$wsResults.Range("B16").Value = "Волна"
$wsResults.Range("C16").Value = "едет"
$wsResults.Range("B17").Value = "Прибой"
$wsResults.Range("C17").Value = "не едет"
$wsResults.Range("B18").Value = "Лесной"
$wsResults.Range("C18").Value = "не едет"
$wsResults.Range("B19").Value = "Лесной"
$wsResults.Range("B20").Value = "Волна"
$wsResults.Range("B21").Value = "Волна"
$wsResults.Range("B22").Value = "Прибой"
$wsResults.Range("B23").Value = "Прибой"
$wsResults.Range("C23").Value = "не едет"
$wsResults.Range("B24").Value = "Искра"
$wsResults.Range("B25").Value = "Янтарь"
$wsResults.Range("B26").Value = "Лесной"
$wsResults.Range("B27").Value = "Лесной"
$wsResults.Range("C27").Value = "не едет"
$wsResults.Range("B28").Value = "Лесной"
$wsResults.Range("B29").Value = "Искра"
$wsResults.Range("B30").Value = "Прибой"
